# Auto-generated Excel COM-interop script
# Applies the "chore: update Sheets via scheduled runner" data refresh
# to columns H:N (price/profit columns) across all 8 sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2932.5334
$ws.Range("I17").Value = 1000
$ws.Range("J17").Value = 3415.6667
$ws.Range("K17").Value = 3000
$ws.Range("L17").Value = 10247.0001
$ws.Range("M17").Value = -2832
$ws.Range("N17").Value = -10583.0001

$ws.Range("H76").Value = 1933
$ws.Range("I76").Value = 1399.5
$ws.Range("J76").Value = 3000
$ws.Range("K76").Value = 1399.5
$ws.Range("L76").Value = 3000
$ws.Range("M76").Value = -1084.5
$ws.Range("N76").Value = -3630

$ws.Range("H79").Value = 1933
$ws.Range("I79").Value = 1399.5
$ws.Range("J79").Value = 3000
$ws.Range("K79").Value = 1399.5
$ws.Range("L79").Value = 3000
$ws.Range("M79").Value = -307.5
$ws.Range("N79").Value = -5184

$ws.Range("H99").Value = 499
$ws.Range("I99").Value = 499
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 1497
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 1
$ws.Range("N99").ClearContents()

$ws.Range("H107").Value = 1999.2858
$ws.Range("I107").Value = 2488.6365
$ws.Range("J107").Value = 205
$ws.Range("K107").Value = 2488.6365
$ws.Range("L107").Value = 205
$ws.Range("M107").Value = -568.6365000000001
$ws.Range("N107").Value = -4045

$ws.Range("H111").Value = 1786.2
$ws.Range("I111").Value = 885
$ws.Range("J111").Value = 5391
$ws.Range("K111").Value = 2655
$ws.Range("L111").Value = 16173
$ws.Range("M111").Value = 412
$ws.Range("N111").Value = -22307

$ws.Range("H113").Value = 7482.7144
$ws.Range("I113").Value = 7666.5
$ws.Range("J113").Value = 6380
$ws.Range("K113").Value = 7666.5
$ws.Range("L113").Value = 6380
$ws.Range("M113").Value = -4412.5
$ws.Range("N113").Value = -12888

$ws.Range("H115").Value = 928.4286
$ws.Range("I115").Value = 928.4286
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 2785.2858
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -1218.2858

$ws.Range("H132").Value = 1182.9474
$ws.Range("I132").Value = 1182.9474
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3548.8422
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1018.8422

$ws.Range("H137").Value = 4749.5
$ws.Range("I137").Value = 4299.4
$ws.Range("J137").Value = 5499.6665
$ws.Range("K137").Value = 12898.2
$ws.Range("L137").Value = 16498.9995
$ws.Range("M137").Value = -10348.2
$ws.Range("N137").Value = -21598.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 4370
$ws.Range("I3").Value = 1555
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 1555
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = -1440
$ws.Range("N3").Value = -10230

$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H74").Value = 21609.262
$ws.Range("I74").Value = 20773.363
$ws.Range("J74").Value = 39999
$ws.Range("K74").Value = 20773.363
$ws.Range("L74").Value = 39999
$ws.Range("M74").Value = -19899.363
$ws.Range("N74").Value = -41747

$ws.Range("H77").Value = 21609.262
$ws.Range("I77").Value = 20773.363
$ws.Range("J77").Value = 39999
$ws.Range("K77").Value = 103866.815
$ws.Range("L77").Value = 199995
$ws.Range("M77").Value = -99498.815
$ws.Range("N77").Value = -208731

$ws.Range("H110").Value = 5576.375
$ws.Range("I110").Value = 2801.5715
$ws.Range("J110").Value = 25000
$ws.Range("K110").Value = 2801.5715
$ws.Range("L110").Value = 25000
$ws.Range("M110").Value = -756.5715
$ws.Range("N110").Value = -29090

$ws.Range("H132").Value = 2185
$ws.Range("I132").Value = 1732.8235
$ws.Range("J132").Value = 3466.1667
$ws.Range("K132").Value = 5198.470499999999
$ws.Range("L132").Value = 10398.5001
$ws.Range("M132").Value = -2668.470499999999
$ws.Range("N132").Value = -15458.5001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 14133
$ws.Range("I20").Value = 2400
$ws.Range("J20").Value = 19999.5
$ws.Range("K20").Value = 2400
$ws.Range("L20").Value = 19999.5
$ws.Range("M20").Value = -2153
$ws.Range("N20").Value = -20493.5

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1818.1666
$ws.Range("I16").Value = 1781.8
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1781.8
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -1494.8
$ws.Range("N16").Value = -2574

$ws.Range("H35").Value = 10796
$ws.Range("I35").Value = 995
$ws.Range("J35").Value = 50000
$ws.Range("K35").Value = 995
$ws.Range("L35").Value = 50000
$ws.Range("M35").Value = -701
$ws.Range("N35").Value = -50588

$ws.Range("H86").Value = 6997.75
$ws.Range("I86").Value = 6997.75
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 6997.75
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -5874.75
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 6997.75
$ws.Range("I89").Value = 6997.75
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 34988.75
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -29372.75
$ws.Range("N89").ClearContents()

$ws.Range("H113").Value = 1818.1666
$ws.Range("I113").Value = 1781.8
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1781.8
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 388.2
$ws.Range("N113").Value = -6340

$ws.Range("H132").Value = 1847.9333
$ws.Range("I132").Value = 1222.5
$ws.Range("J132").Value = 3098.8
$ws.Range("K132").Value = 3667.5
$ws.Range("L132").Value = 9296.400000000001
$ws.Range("M132").Value = -1137.5
$ws.Range("N132").Value = -14356.4

$ws.Range("H134").Value = 5434.5
$ws.Range("I134").Value = 4825.1113
$ws.Range("J134").Value = 7262.6665
$ws.Range("K134").Value = 14475.3339
$ws.Range("L134").Value = 21787.9995
$ws.Range("M134").Value = -11940.3339
$ws.Range("N134").Value = -26857.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 212.5
$ws.Range("I26").Value = 266.66666
$ws.Range("J26").Value = 50
$ws.Range("K26").Value = 799.9999799999999
$ws.Range("L26").Value = 150
$ws.Range("M26").Value = -511.9999799999999
$ws.Range("N26").Value = -726

$ws.Range("H80").Value = 10665.667
$ws.Range("I80").Value = 10665.667
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 31997.001
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -31061.001

$ws.Range("H83").Value = 10665.667
$ws.Range("I83").Value = 10665.667
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 95991.003
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -91311.003

$ws.Range("H86").Value = 1130.75
$ws.Range("I86").Value = 1414
$ws.Range("J86").Value = 1090.2858
$ws.Range("K86").Value = 4242
$ws.Range("L86").Value = 3270.8574
$ws.Range("M86").Value = -3056
$ws.Range("N86").Value = -5642.857400000001

$ws.Range("H89").Value = 1130.75
$ws.Range("I89").Value = 1414
$ws.Range("J89").Value = 1090.2858
$ws.Range("K89").Value = 12726
$ws.Range("L89").Value = 9812.572200000001
$ws.Range("M89").Value = -6798
$ws.Range("N89").Value = -21668.5722

$ws.Range("H122").Value = 3074.3076
$ws.Range("I122").Value = 912.5
$ws.Range("J122").Value = 3467.3635
$ws.Range("K122").Value = 8212.5
$ws.Range("L122").Value = 31206.2715
$ws.Range("M122").Value = -5762.5
$ws.Range("N122").Value = -36106.2715

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 317.66666
$ws.Range("I2").Value = 501.8
$ws.Range("J2").Value = 87.5
$ws.Range("K2").Value = 501.8
$ws.Range("L2").Value = 87.5
$ws.Range("M2").Value = -388.8
$ws.Range("N2").Value = -313.5

$ws.Range("H102").Value = 2528.6667
$ws.Range("I102").Value = 2528.6667
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2528.6667
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -906.6667000000002

$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H126").Value = 6266.6665
$ws.Range("I126").Value = 6266.6665
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 18799.9995
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -16329.9995

$ws.Range("H132").Value = 1651.8948
$ws.Range("I132").Value = 1303.375
$ws.Range("J132").Value = 3510.6667
$ws.Range("K132").Value = 3910.125
$ws.Range("L132").Value = 10532.0001
$ws.Range("M132").Value = -1380.125
$ws.Range("N132").Value = -15592.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2000
$ws.Range("I16").Value = 2000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1830

$ws.Range("H18").Value = 54997
$ws.Range("I18").Value = 9995
$ws.Range("J18").Value = 99999
$ws.Range("K18").Value = 9995
$ws.Range("L18").Value = 99999
$ws.Range("M18").Value = -9823
$ws.Range("N18").Value = -100343

$ws.Range("H32").Value = 1406.5
$ws.Range("I32").Value = 13
$ws.Range("J32").Value = 2800
$ws.Range("K32").Value = 13
$ws.Range("L32").Value = 2800
$ws.Range("M32").Value = 304
$ws.Range("N32").Value = -3434

$ws.Range("H100").Value = 4666.6665
$ws.Range("I100").Value = 4666.6665
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 4666.6665
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -4125.6665

$ws.Range("H136").Value = 3998.6667
$ws.Range("I136").Value = 4009.923
$ws.Range("J136").Value = 3925.5
$ws.Range("K136").Value = 12029.769
$ws.Range("L136").Value = 11776.5
$ws.Range("M136").Value = -9479.769
$ws.Range("N136").Value = -16876.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 5000000
$ws.Range("I3").Value = 5000000
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5000000
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -4999886

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H107").Value = 840.0714
$ws.Range("I107").Value = 621.5
$ws.Range("J107").Value = 1004
$ws.Range("K107").Value = 1864.5
$ws.Range("L107").Value = 3012
$ws.Range("M107").Value = 55.5
$ws.Range("N107").Value = -6852

$ws.Range("H136").Value = 1497.9678
$ws.Range("I136").Value = 1497.9678
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4493.903399999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -1943.903399999999
